$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New catalog rows (MCH155-1 .. MCH155-6) describing the "AMANDLA" series.
# Columns: A=identifier, C=title, D=date_s, E=levelOfDescription,
#          F=extentAndMedium, G=notes, H=file_path (left blank)
# ---------------------------------------------------------------------------

$rows = @(
    @{ id = "MCH155-1"; title = "ANGOLA BULLETIN, AMANDLA 1971-1979"; date = "1971-1979" },
    @{ id = "MCH155-2"; title = "AMANDLA (BOA, TAIRAS) 1980-1983"; date = "1980-1983" },
    @{ id = "MCH155-3"; title = "AMANDLA ( BOA, KZA, KAIROS) 1984-1987"; date = "1984-1987" },
    @{ id = "MCH155-4"; title = "AMANDLA (BOA,KZA,KAIROS) 1988-1991"; date = "1988-1991" },
    @{ id = "MCH155-5"; title = "AMANDLA (KZA, KAIROS) 1991 BOMEED, 1992 UNBOUNDED"; date = "1991" },
    @{ id = "MCH155-6"; title = "VARIOUS (3RD) DUPLICATES OF AMANDLA"; date = "" }
)

$level = "Series"
$extent = "1 Box"
$notes = "LOCATION: 21O | GRAP COUNT NUMER: NONE"

$r = 2
foreach ($row in $rows) {

    $ws.Cells.Item($r, 1).Value = $row.id
    $ws.Cells.Item($r, 1).Font.ThemeColor = 1
    $ws.Cells.Item($r, 1).Font.Name = "Calibri"
    $ws.Cells.Item($r, 1).Font.Size = 10

    $ws.Cells.Item($r, 3).Value = $row.title
    $ws.Cells.Item($r, 3).Font.ThemeColor = 1
    $ws.Cells.Item($r, 3).Font.Name = "Calibri"
    $ws.Cells.Item($r, 3).Font.Size = 10

    if ($row.date -ne "") {
        # Values that look purely numeric (e.g. "1991") would otherwise be
        # coerced to a number; round-trip them through a text formula so
        # they land back in the sheet as shared-string text, matching the
        # source data's string typing, without leaving a stray NumberFormat
        # or quote-prefix behind.
        $looksNumeric = $row.date -match '^-?[0-9]+(\.[0-9]+)?$'

        if ($looksNumeric) {
            $ws.Cells.Item($r, 4).Formula = '="' + $row.date + '"'
            $ws.Cells.Item($r, 4).Copy()
            $ws.Cells.Item($r, 4).PasteSpecial(-4163)
        } else {
            $ws.Cells.Item($r, 4).Value = $row.date
        }
    }
    $ws.Cells.Item($r, 4).Font.ThemeColor = 1
    $ws.Cells.Item($r, 4).Font.Name = "Calibri"
    $ws.Cells.Item($r, 4).Font.Size = 10

    $ws.Cells.Item($r, 5).Value = $level
    $ws.Cells.Item($r, 5).Font.ThemeColor = 1
    $ws.Cells.Item($r, 5).Font.Name = "Calibri"
    $ws.Cells.Item($r, 5).Font.Size = 10

    $ws.Cells.Item($r, 6).Value = $extent
    $ws.Cells.Item($r, 6).Font.ThemeColor = 1
    $ws.Cells.Item($r, 6).Font.Name = "Calibri"
    $ws.Cells.Item($r, 6).Font.Size = 10
    $ws.Cells.Item($r, 6).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 7).Value = $notes
    $ws.Cells.Item($r, 7).Font.ThemeColor = 1
    $ws.Cells.Item($r, 7).Font.Name = "Calibri"
    $ws.Cells.Item($r, 7).Font.Size = 10

    # column H (file_path) is left blank but still carries the row format
    $ws.Cells.Item($r, 8).Font.ThemeColor = 1
    $ws.Cells.Item($r, 8).Font.Name = "Calibri"
    $ws.Cells.Item($r, 8).Font.Size = 10

    $ws.Rows.Item($r).RowHeight = 15.75

    $r = $r + 1
}

# Mirror the author's final selection over the newly entered block.
$ws.Range("A2:K7").Select()
